# feat: add 2022-Q1 data
#
# 1. Insert a brand-new sheet "2022-Q1" between "2021-Q4" and "总计",
#    populated with the fund holdings data for the new quarter.
# 2. Update the "总计" (totals) sheet to add a new row for "2022-Q1"
#    above the existing "2021-Q4" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Re-create the "总计" sheet so that, after we insert "2022-Q1" before
# it, the sheetId ordering matches: 2021-Q4=1, 2022-Q1=2, 总计=3.
# ---------------------------------------------------------------------
$wsTotalOld = $wb.Worksheets.Item("总计")
[void]$wsTotalOld.Delete()

$wsQ4 = $wb.Worksheets.Item(1)

$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

# Match the original page-margin settings used on these two sheets
# (left/right = 0.75in, top/bottom = 1in, header/footer = 0.5in).
foreach ($sheet in @($wsQ1, $wsTotal)) {
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
}

# ---------------------------------------------------------------------
# Populate the new "2022-Q1" sheet - same layout as "2021-Q4".
# ---------------------------------------------------------------------
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$col = 2
foreach ($h in $headers) {
    $cell = $wsQ1.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $col++
}

$rows = @(
    @("013233", "华夏中证500指数智选增强A", "39.48", "92.73", "1.36", "0.5369", 8),
    @("007994", "华夏中证500指数增强A",     "31.45", "92.72", "1.35", "0.4246", 8),
    @("007995", "华夏中证500指数增强C",     "5.45",  "92.72", "1.35", "0.0736", 8),
    @("013234", "华夏中证500指数智选增强C", "4.28",  "92.73", "1.36", "0.0582", 8),
    @("501219", "华夏智胜先锋股票（LOF）A", "3.61",  "94.50", "1.40", "0.0505", 1),
    @("014198", "华夏智胜先锋股票（LOF）C", "1.30",  "94.50", "1.40", "0.0182", 1)
)

$r = 2
$idx = 0
foreach ($row in $rows) {
    $acell = $wsQ1.Cells.Item($r, 1)
    $acell.Value = $idx
    $acell.Font.Bold = $true
    $acell.HorizontalAlignment = -4108
    $acell.VerticalAlignment = -4160
    $acell.Borders.LineStyle = 1

    $textRange = $wsQ1.Range($wsQ1.Cells.Item($r, 2), $wsQ1.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"
    $wsQ1.Cells.Item($r, 2).Value = $row[0]
    $wsQ1.Cells.Item($r, 3).Value = $row[1]
    $wsQ1.Cells.Item($r, 4).Value = $row[2]
    $wsQ1.Cells.Item($r, 5).Value = $row[3]
    $wsQ1.Cells.Item($r, 6).Value = $row[4]
    $wsQ1.Cells.Item($r, 7).Value = $row[5]
    $textRange.ClearFormats()

    $wsQ1.Cells.Item($r, 8).Value = $row[6]

    $r++
    $idx++
}

# ---------------------------------------------------------------------
# Populate the "总计" sheet.
# ---------------------------------------------------------------------
$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
$col = 2
foreach ($h in $totalHeaders) {
    $cell = $wsTotal.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $col++
}

$totalRows = @(
    @("2022-Q1", 6, 1.16),
    @("2021-Q4", 2, 0.09)
)

$r = 2
$idx = 0
foreach ($row in $totalRows) {
    $acell = $wsTotal.Cells.Item($r, 1)
    $acell.Value = $idx
    $acell.Font.Bold = $true
    $acell.HorizontalAlignment = -4108
    $acell.VerticalAlignment = -4160
    $acell.Borders.LineStyle = 1

    $wsTotal.Cells.Item($r, 2).Value = $row[0]
    $wsTotal.Cells.Item($r, 3).Value = $row[1]
    $wsTotal.Cells.Item($r, 4).Value = $row[2]

    $r++
    $idx++
}

# Restore "2021-Q4" as the active/selected sheet, matching the original.
[void]$wsQ4.Select()

